$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.134.14'
$ws.Range("E2").Value = '  +2.84%  '
$ws.Range("D3").Value = '2.619.47'
$ws.Range("E3").Value = '  +1.00%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.17'
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.73'
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").Value = '2.642.43'
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.57'
$ws.Range("E10").Value = '  -1.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.107'
$ws.Range("E11").Value = '  +2.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.157'
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.373'
$ws.Range("E13").Value = '  +7.33%  '
$ws.Range("D14").Value = '3.085.94'
$ws.Range("E14").Value = '  +1.08%  '
$ws.Range("D15").Value = '61.125.73'
$ws.Range("E15").Value = '  +2.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.59'
$ws.Range("E16").Value = '  +4.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000142'
$ws.Range("E17").Value = '  +3.07%  '
$ws.Range("D18").Value = '2.632.56'
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.73'
$ws.Range("E19").Value = '  +3.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.16'
$ws.Range("E20").Value = '  +8.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '353.83'
$ws.Range("E21").Value = '  +4.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.07'
$ws.Range("E22").Value = '  +13.42%  '
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("E24").Value = '  +12.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.93'
$ws.Range("E25").Value = '  -1.03%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.162'
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.79'
$ws.Range("E28").Value = '  +6.60%  '
$ws.Range("D29").Value = '0.0₃0803'
$ws.Range("E29").Value = '  +2.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.86'
$ws.Range("E30").Value = '  +10.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.31'
$ws.Range("E32").Value = '  +3.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '161.82'
$ws.Range("E33").Value = '  +1.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.55'
$ws.Range("E34").Value = '  +2.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.27'
$ws.Range("E35").Value = '  +5.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.970'
$ws.Range("E36").Value = '  +9.31%  '
$ws.Range("E37").Value = '  +4.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.61'
$ws.Range("E38").Value = '  +7.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.76'
$ws.Range("E39").Value = '  +1.58%  '
$ws.Range("E40").Value = '  -2.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.82'
$ws.Range("E41").Value = '  +3.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '302.22'
$ws.Range("E42").Value = '  +1.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '140.95'
$ws.Range("E43").Value = '  +12.50%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0989'
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.994'
$ws.Range("E45").Value = '  -0.43%  '
$ws.Range("E46").Value = '  +1.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0552'
$ws.Range("E47").Value = '  +2.30%  '
$ws.Range("E48").Value = '  +3.68%  '
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.72'
$ws.Range("E50").Value = '  +5.66%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.054.47'
$ws.Range("E51").Value = '  +5.16%  '
